# Apply the "New crime data collected" weekly update to the 63rd Precinct sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text: "Volume 31   Number  22" -> "...  23"
#              "Report Covering the Week 5/27/2024 Through 6/2/2024"
#           -> "Report Covering the Week 6/3/2024  Through 6/9/2024"
# ---------------------------------------------------------------------
$a8 = $ws.Range("A8")
$a8full = $a8.Value2
$a8.Characters($a8full.Length - 1, 2).Text = "23"

$c9 = $ws.Range("C9")
$c9.Characters(27, 9).Text = "6/3/2024"
$c9.Characters(46, 8).Text = "6/9/2024"

# ---------------------------------------------------------------------
# A few cells flip between a numeric count and the "no data" text
# placeholders ("0" / "***.*") used throughout this sheet. Setting
# .Value with a leading apostrophe forces text, then PasteSpecial of
# the formats (only) from an existing same-style cell keeps the cell
# style identical to the rest of the text-placeholder cells / number
# cells in this table.
# ---------------------------------------------------------------------
function Set-TextCell($cell, $text, $donor) {
    $cell.Value = "'" + $text
    $donor.Copy()
    $cell.PasteSpecial(-4122)
}

$textStyleDonor = $ws.Range("D15")     # "0"-style text placeholder cell
$starStyleDonor = $ws.Range("E15")     # "***.*"-style text placeholder cell
$numStyleDonor  = $ws.Range("D16")     # plain numeric cell

Set-TextCell $ws.Range("C15") "0" $textStyleDonor
Set-TextCell $ws.Range("C18") "0" $textStyleDonor
Set-TextCell $ws.Range("D23") "0" $textStyleDonor
Set-TextCell $ws.Range("E23") "***.*" $starStyleDonor

$ws.Range("C23").Value = 1
$numStyleDonor.Copy()
$ws.Range("C23").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Remaining plain numeric value updates (counts and computed % changes)
# ---------------------------------------------------------------------
$ws.Range("N15").Value = -66.666666666666
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -75
$ws.Range("I16").Value = 61
$ws.Range("J16").Value = 45
$ws.Range("K16").Value = 35.555555555555
$ws.Range("L16").Value = 24.489795918367
$ws.Range("M16").Value = -37.755102040816
$ws.Range("N16").Value = -78.892733564013
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = -50
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = 36.363636363636
$ws.Range("I17").Value = 66
$ws.Range("J17").Value = 67
$ws.Range("K17").Value = -1.492537313432
$ws.Range("L17").Value = -5.714285714285
$ws.Range("M17").Value = 17.857142857142
$ws.Range("N17").Value = -51.111111111111
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -100
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = -33.333333333333
$ws.Range("J18").Value = 46
$ws.Range("K18").Value = -39.130434782608
$ws.Range("L18").Value = -37.777777777777
$ws.Range("N18").Value = -94.954954954955
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 47
$ws.Range("G19").Value = 48
$ws.Range("H19").Value = -2.083333333333
$ws.Range("I19").Value = 262
$ws.Range("J19").Value = 268
$ws.Range("K19").Value = -2.238805970149
$ws.Range("L19").Value = 26.570048309178
$ws.Range("M19").Value = 20.73732718894
$ws.Range("N19").Value = -12.956810631229
$ws.Range("F20").Value = 14
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = 55.555555555555
$ws.Range("I20").Value = 68
$ws.Range("J20").Value = 45
$ws.Range("K20").Value = 51.111111111111
$ws.Range("L20").Value = 36
$ws.Range("M20").Value = -6.849315068493
$ws.Range("N20").Value = -94.598888006354
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 23
$ws.Range("E21").Value = -21.739130434782
$ws.Range("F21").Value = 90
$ws.Range("G21").Value = 88
$ws.Range("H21").Value = 2.272727272727
$ws.Range("I21").Value = 489
$ws.Range("J21").Value = 477
$ws.Range("K21").Value = 2.51572327044
$ws.Range("L21").Value = 15.058823529411
$ws.Range("M21").Value = -12.208258527827
$ws.Range("N21").Value = -80.905896134322
$ws.Range("F23").Value = 3
$ws.Range("H23").Value = 200
$ws.Range("I23").Value = 17
$ws.Range("K23").Value = 112.5
$ws.Range("L23").Value = 21.428571428571
$ws.Range("M23").Value = 13.333333333333
$ws.Range("C24").Value = 25
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = 13.636363636363
$ws.Range("F24").Value = 123
$ws.Range("H24").Value = 55.696202531645
$ws.Range("I24").Value = 673
$ws.Range("J24").Value = 520
$ws.Range("K24").Value = 29.423076923076
$ws.Range("L24").Value = 47.587719298245
$ws.Range("M24").Value = 69.521410579345
$ws.Range("D25").Value = 19
$ws.Range("E25").Value = 5.263157894736
$ws.Range("G25").Value = 68
$ws.Range("H25").Value = 52.941176470588
$ws.Range("I25").Value = 560
$ws.Range("J25").Value = 349
$ws.Range("K25").Value = 60.458452722063
$ws.Range("L25").Value = 91.126279863481
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 900
$ws.Range("F26").Value = 34
$ws.Range("G26").Value = 22
$ws.Range("H26").Value = 54.545454545454
$ws.Range("I26").Value = 141
$ws.Range("J26").Value = 114
$ws.Range("K26").Value = 23.684210526315
$ws.Range("L26").Value = 46.875
$ws.Range("M26").Value = -1.398601398601
$ws.Range("C27").Value = 2
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = 200
$ws.Range("I27").Value = 6
$ws.Range("K27").Value = -45.454545454545
$ws.Range("L27").Value = -33.333333333333
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 200
$ws.Range("L28").Value = -6.666666666666
$ws.Range("N29").Value = -91.304347826087
$ws.Range("N30").Value = -90.47619047619
$ws.Range("C31").Value = 1
$ws.Range("F31").Value = 3
$ws.Range("I31").Value = 6
$ws.Range("K31").Value = 500
$ws.Range("L31").Value = 20
